$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for three new columns (R,S,T) before the existing lat/lon
#     data columns (old R,S), which shift right to become U,V ---
[void]$ws.Range("R1:T1").EntireColumn.Insert()

# --- New header row labels for the inserted columns ---
$ws.Range("R1").Value = "Photograph"
$ws.Range("S1").Value = "Letter of Interest or Intent"
$ws.Range("T1").Value = "Pathology graduates only"

# --- Row 41: Baylor College of Medicine ---
$ws.Range("A41").Value = "Baylor College of Medicine"
$ws.Range("B41").Value = "Houston"
$ws.Range("C41").Value = "Texas"
$ws.Range("D41").Value = "West South Central"
$ws.Range("E41").Value = $false
$ws.Range("F41").Value = $true
$ws.Range("G41").Value = $true
$ws.Range("H41").Value = $true
$ws.Range("I41").Value = 2
$ws.Range("J41").Value = $true
$ws.Range("K41").Value = $true
$ws.Range("R41").Value = $true

# --- Row 42: University at Buffalo ---
$ws.Range("A42").Value = "University at Buffalo"
$ws.Range("B42").Value = "Buffalo"
$ws.Range("C42").Value = "New York"
$ws.Range("D42").Value = "Mid Atlantic"
$ws.Range("E42").Value = $false
$ws.Range("F42").Value = $true
$ws.Range("G42").Value = $true
$ws.Range("H42").Value = $true
$ws.Range("I42").Value = 3

# --- Row 43: McGaw Medical Center of Northwestern University Program ---
$ws.Range("A43").Value = "McGaw Medical Center of Northwestern University Program"
$ws.Range("B43").Value = "Chicago"
$ws.Range("C43").Value = "Illinois"
$ws.Range("D43").Value = "East North Central"
$ws.Range("E43").Value = $false

# --- Row 44: UC Davis ---
$ws.Range("A44").Value = "UC Davis"
$ws.Range("B44").Value = "Sacramento"
$ws.Range("C44").Value = "California"
$ws.Range("D44").Value = "Pacific"
$ws.Range("E44").Value = $false
$ws.Range("F44").Value = $true
$ws.Range("G44").Value = $true
$ws.Range("I44").Value = 3
$ws.Range("L44").Value = $true
$ws.Range("S44").Value = $true

# --- Row 45: Massachusettes General Hospital ---
$ws.Range("A45").Value = "Massachusettes General Hospital"
$ws.Range("B45").Value = "Boston"
$ws.Range("C45").Value = "Massachusetts"
$ws.Range("D45").Value = "New England"
$ws.Range("E45").Value = $false
$ws.Range("F45").Value = $true
$ws.Range("G45").Value = $true
$ws.Range("I45").Value = 3
$ws.Range("K45").Value = $true
$ws.Range("L45").Value = $true
$ws.Range("T45").Value = "Yes"

# --- Final selection, matching the state left behind after the edits ---
[void]$ws.Range("T46").Select()
